# modificato scarico per calcolare valori batch
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "numero serventi" changes from 4 to 10
$ws.Range("B2").Value = 10

$excel.CalculateFullRebuild()
